$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3541.9473
$ws.Range("I62").Value = 3257.6667
$ws.Range("J62").Value = 4029.2856
$ws.Range("K62").Value = 3257.6667
$ws.Range("L62").Value = 4029.2856
$ws.Range("M62").Value = -2633.6667
$ws.Range("N62").Value = -5277.2856

$ws.Range("H65").Value = 3541.9473
$ws.Range("I65").Value = 3257.6667
$ws.Range("J65").Value = 4029.2856
$ws.Range("K65").Value = 16288.3335
$ws.Range("L65").Value = 20146.428
$ws.Range("M65").Value = -13168.3335
$ws.Range("N65").Value = -26386.428

$ws.Range("H99").Value = 436.63635
$ws.Range("I99").Value = 323.83334
$ws.Range("J99").Value = 572
$ws.Range("K99").Value = 971.5000200000001
$ws.Range("L99").Value = 1716
$ws.Range("M99").Value = 526.4999799999999
$ws.Range("N99").Value = -4712

$ws.Range("H125").Value = 4017.125
$ws.Range("I125").Value = 2277
$ws.Range("J125").Value = 5061.2
$ws.Range("K125").Value = 20493
$ws.Range("L125").Value = 45550.8
$ws.Range("M125").Value = -18033
$ws.Range("N125").Value = -50470.8

$ws.Range("H129").Value = 1105.72
$ws.Range("I129").Value = 598.5
$ws.Range("J129").Value = 1116.0714
$ws.Range("K129").Value = 1795.5
$ws.Range("L129").Value = 3348.2142
$ws.Range("M129").Value = 3204.5
$ws.Range("N129").Value = -13348.2142

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13336915
$ws.Range("I32").Value = 14495081
$ws.Range("K32").Value = 14495081
$ws.Range("M32").Value = -14494794

$ws.Range("H45").Value = 1267.091
$ws.Range("I45").Value = 896.1818
$ws.Range("J45").Value = 1638
$ws.Range("K45").Value = 896.1818
$ws.Range("L45").Value = 1638
$ws.Range("M45").Value = -519.1818
$ws.Range("N45").Value = -2392

$ws.Range("H96").Value = 23668.8
$ws.Range("J96").Value = 23668.8
$ws.Range("L96").Value = 23668.8
$ws.Range("N96").Value = -29160.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 39800
$ws.Range("J16").Value = 39800
$ws.Range("L16").Value = 39800
$ws.Range("N16").Value = -40140

$ws.Range("H134").Value = 1691.697
$ws.Range("I134").Value = 1719.5186
$ws.Range("K134").Value = 5158.5558
$ws.Range("M134").Value = -2623.5558

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 7502125
$ws.Range("I6").Value = 30000000
$ws.Range("J6").Value = 2833.3333
$ws.Range("K6").Value = 30000000
$ws.Range("L6").Value = 2833.3333
$ws.Range("M6").Value = -29999887
$ws.Range("N6").Value = -3059.3333

$ws.Range("H7").Value = 35.636364
$ws.Range("I7").Value = 22.833334
$ws.Range("J7").Value = 51
$ws.Range("K7").Value = 22.833334
$ws.Range("L7").Value = 51
$ws.Range("M7").Value = 90.16666599999999
$ws.Range("N7").Value = -277

$ws.Range("H17").Value = 10166.667
$ws.Range("I17").Value = 10000
$ws.Range("J17").Value = 10250
$ws.Range("K17").Value = 10000
$ws.Range("L17").Value = 10250
$ws.Range("M17").Value = -9826
$ws.Range("N17").Value = -10598

$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()

$ws.Range("H31").Value = 1894.1177
$ws.Range("I31").Value = 1353.95
$ws.Range("J31").Value = 3858.3635
$ws.Range("K31").Value = 1353.95
$ws.Range("L31").Value = 3858.3635
$ws.Range("M31").Value = -1058.95
$ws.Range("N31").Value = -4448.363499999999

$ws.Range("H34").Value = 1894.1177
$ws.Range("I34").Value = 1353.95
$ws.Range("J34").Value = 3858.3635
$ws.Range("K34").Value = 1353.95
$ws.Range("L34").Value = 3858.3635
$ws.Range("M34").Value = -1151.95
$ws.Range("N34").Value = -4262.363499999999

$ws.Range("H41").Value = 9949.166999999999
$ws.Range("I41").Value = 8000
$ws.Range("J41").Value = 10339
$ws.Range("K41").Value = 8000
$ws.Range("L41").Value = 10339
$ws.Range("M41").Value = -7572
$ws.Range("N41").Value = -11195

$ws.Range("H50").Value = 19932.666
$ws.Range("J50").Value = 19932.666
$ws.Range("L50").Value = 19932.666
$ws.Range("N50").Value = -21182.666

$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("N51").ClearContents()

$ws.Range("H58").Value = 832.5
$ws.Range("I58").Value = 712.6531
$ws.Range("J58").Value = 2007
$ws.Range("K58").Value = 712.6531
$ws.Range("L58").Value = 2007
$ws.Range("M58").Value = -509.6531
$ws.Range("N58").Value = -2413

$ws.Range("H59").Value = 23063.5
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 23063.5
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 23063.5
$ws.Range("M59").ClearContents()
$ws.Range("N59").Value = -25353.5

$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").ClearContents()

$ws.Range("H62").Value = 2636.9092
$ws.Range("I62").Value = 2883.3333
$ws.Range("J62").Value = 2341.2
$ws.Range("K62").Value = 2883.3333
$ws.Range("L62").Value = 2341.2
$ws.Range("M62").Value = -2259.3333
$ws.Range("N62").Value = -3589.2

$ws.Range("H65").Value = 2636.9092
$ws.Range("I65").Value = 2883.3333
$ws.Range("J65").Value = 2341.2
$ws.Range("K65").Value = 14416.6665
$ws.Range("L65").Value = 11706
$ws.Range("M65").Value = -11296.6665
$ws.Range("N65").Value = -17946

$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws.Range("H74").Value = 18595.555
$ws.Range("J74").Value = 19670
$ws.Range("L74").Value = 19670
$ws.Range("N74").Value = -21418

$ws.Range("H77").Value = 18595.555
$ws.Range("J77").Value = 19670
$ws.Range("L77").Value = 59010
$ws.Range("N77").Value = -67746

$ws.Range("H136").Value = 832.5
$ws.Range("I136").Value = 712.6531
$ws.Range("J136").Value = 2007
$ws.Range("K136").Value = 2137.9593
$ws.Range("L136").Value = 6021
$ws.Range("M136").Value = 412.0407
$ws.Range("N136").Value = -11121

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 762.7778
$ws.Range("I131").Value = 436.69232
$ws.Range("J131").Value = 895.25
$ws.Range("K131").Value = 1310.07696
$ws.Range("L131").Value = 2685.75
$ws.Range("M131").Value = 3729.92304
$ws.Range("N131").Value = -12765.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H128").Value = 27000
$ws.Range("J128").Value = 27000
$ws.Range("L128").Value = 27000
$ws.Range("N128").Value = -36960

$ws.Range("H132").Value = 2666.1592
$ws.Range("I132").Value = 2449.258
$ws.Range("J132").Value = 3183.3845
$ws.Range("K132").Value = 7347.773999999999
$ws.Range("L132").Value = 9550.1535
$ws.Range("M132").Value = -4817.773999999999
$ws.Range("N132").Value = -14610.1535

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2150
$ws.Range("I81").Value = 1000
$ws.Range("J81").Value = 2533.3333
$ws.Range("K81").Value = 2000
$ws.Range("L81").Value = 5066.6666
$ws.Range("M81").Value = -939
$ws.Range("N81").Value = -7188.6666

$ws.Range("H84").Value = 2150
$ws.Range("I84").Value = 1000
$ws.Range("J84").Value = 2533.3333
$ws.Range("K84").Value = 10000
$ws.Range("L84").Value = 25333.333
$ws.Range("M84").Value = -4696
$ws.Range("N84").Value = -35941.333
